# EWD-22219 - Build Experience: Build/Download Experience MarkUp
#
# Capitalize the EN/NL/DE labels used for the experience build status
# (Building…/Failed/Complete) and re-order the shared-string entries so
# the "building" row also carries an English string (it was previously
# missing one), matching the other rows in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 31: experienceBuildingStatus -> Building… / Samenstellen… / Aufbau…
$ws.Range("C31").Value = "Building…"
$ws.Range("D31").Value = "Samenstellen…"
$ws.Range("E31").Value = "Aufbau…"

# Row 32: experienceFailedStatus -> Failed / Mislukt / Fehlgeschlagen
$ws.Range("C32").Value = "Failed"
$ws.Range("D32").Value = "Mislukt"
$ws.Range("E32").Value = "Fehlgeschlagen"

# Row 33: experienceCompleteStatus -> Complete / Voltooid / Vollständig
# (entered DE/NL/EN, matching how these values were originally typed)
$ws.Range("E33").Value = "Vollständig"
$ws.Range("D33").Value = "Voltooid"
$ws.Range("C33").Value = "Complete"

# Scroll/selection state recorded by Excel when the sheet was last saved.
$window = $excel.ActiveWindow
$window.ScrollColumn = 2
$window.ScrollRow = 19
$ws.Range("C35").Select()
